$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("CoCRoI")
$ws.Range("B2").Value2 = -7.852283465351582
$ws.Range("C2").Value2 = -3444.207834989838
$ws.Range("B3").Value2 = -1.952164101635395
$ws.Range("C3").Value2 = -856.2679790798252
$ws.Range("B4").Value2 = -2.158498741133831
$ws.Range("C4").Value2 = -946.7715103298265
$ws.Range("B5").Value2 = -2.037905112250963
$ws.Range("C5").Value2 = -893.8761298610789
$ws.Range("B6").Value2 = -1.906826571400646
$ws.Range("C6").Value2 = -836.3818048806086
$ws.Range("B7").Value2 = -1.765030491594487
$ws.Range("C7").Value2 = -774.1864993756319
$ws.Range("B8").Value2 = -1.612279610001771
$ws.Range("C8").Value2 = -707.1861439370266
$ws.Range("B9").Value2 = -1.448331959246771
$ws.Range("C9").Value2 = -635.2746056246151
$ws.Range("B10").Value2 = -1.272940798876131
$ws.Range("C10").Value2 = -558.3436579070431
$ws.Range("B11").Value2 = -1.085854547076172
$ws.Range("C11").Value2 = -476.2829507112858
$ws.Range("B12").Value2 = -0.8868167127249336
$ws.Range("C12").Value2 = -388.979980618974
$ws.Range("B13").Value2 = -0.675565827868966
$ws.Range("C13").Value2 = -296.3200612490252
$ws.Range("B14").Value2 = -0.4518353807204318
$ws.Range("C14").Value2 = -198.1862938684994
$ws.Range("B15").Value2 = -0.2153537492757411
$ws.Range("C15").Value2 = -94.45953827607195
$ws.Range("B16").Value2 = 0.03415586433664432
$ws.Range("C16").Value2 = 14.98161599466062
$ws.Range("B17").Value2 = 0.2969754986655694
$ws.Range("C17").Value2 = 130.2608781021854
$ws.Range("B18").Value2 = 0.5733924967834013
$ws.Range("C18").Value2 = 251.5042839016194
$ws.Range("B19").Value2 = 0.8636995690415287
$ws.Range("C19").Value2 = 378.8402234708405
$ws.Range("B20").Value2 = 1.168194854600343
$ws.Range("C20").Value2 = 512.3994680990755
$ws.Range("B21").Value2 = 1.487181981591563
$ws.Range("C21").Value2 = 652.3151966755995
$ws.Range("B22").Value2 = 1.820970125762113
$ws.Range("C22").Value2 = 798.723021412407
$ws.Range("B23").Value2 = 2.169874067441133
$ws.Range("C23").Value2 = 951.7610128313668
$ws.Range("B24").Value2 = 2.534214246661863
$ws.Range("C24").Value2 = 1111.56972394206
$ws.Range("B25").Value2 = 2.91431681626227
$ws.Range("C25").Value2 = 1278.292213533038
$ws.Range("B26").Value2 = 3.310513692777216
$ws.Range("C26").Value2 = 1452.074068494407
$ws.Range("B27").Value2 = 3.72314260492597
$ws.Range("C27").Value2 = 1633.063425085654
$ws.Range("B28").Value2 = 4.152547139487621
$ws.Range("C28").Value2 = 1821.410989057758
$ws.Range("B29").Value2 = 4.599076784346256
$ws.Range("C29").Value2 = 2017.270054533877
$ws.Range("B30").Value2 = 5.063086968476059
$ws.Range("C30").Value2 = 2220.796521547811
$ws.Range("B31").Value2 = 5.544939098624205
$ws.Range("C31").Value2 = 2432.148912134042
$ws.Range("B32").Value2 = 0.4528956112671862
$ws.Range("C32").Value2 = 198.6513374920696
$ws = $wb.Worksheets.Item("Overall CAGR")
$ws.Range("B2").Value2 = -3397.126729079823
$ws.Range("C2").Value2 = -3397.126729079823
$ws.Range("G2").Value2 = 22331.22580205838
$ws.Range("B3").Value2 = -4388.598458159649
$ws.Range("C3").Value2 = -991.4717290798253
$ws.Range("G3").Value2 = 25165.49861547235
$ws.Range("H3").Value2 = -42.62639244121436
$ws.Range("B4").Value2 = -5331.306312239475
$ws.Range("C4").Value2 = -942.7078540798266
$ws.Range("G4").Value2 = 28187.13242737833
$ws.Range("H4").Value2 = -19.83611481913701
$ws.Range("B5").Value2 = -6220.743194444304
$ws.Range("C5").Value2 = -889.4368822048291
$ws.Range("G5").Value2 = 31407.9769839877
$ws.Range("H5").Value2 = -10.5358868994265
$ws.Range("B6").Value2 = -7052.301810477256
$ws.Range("C6").Value2 = -831.5586160329522
$ws.Range("G6").Value2 = 34840.46559303495
$ws.Range("H6").Value2 = -5.594424899848416
$ws.Range("B7").Value2 = -7821.272662884041
$ws.Range("C7").Value2 = -768.9708524067846
$ws.Range("B8").Value2 = -8522.842015142
$ws.Range("C8").Value2 = -701.5693522579577
$ws.Range("G8").Value2 = 42393.24723876181
$ws.Range("H8").Value2 = -0.5662350551105422
$ws.Range("B9").Value2 = -9152.089825618128
$ws.Range("C9").Value2 = -629.2478104761294
$ws.Range("G9").Value2 = 46541.66563136013
$ws.Range("H9").Value2 = 0.8505713411053994
$ws.Range("B10").Value2 = -9703.987651473435
$ws.Range("C10").Value2 = -551.8978258553079
$ws.Range("G10").Value2 = 50958.06515129079
$ws.Range("H10").Value2 = 1.891966364283193
$ws.Range("B11").Value2 = -10173.39652262617
$ws.Range("C11").Value2 = -469.4088711527384
$ws.Range("G11").Value2 = 55658.39734655086
$ws.Range("H11").Value2 = 2.681696413227197
$ws.Range("B12").Value2 = -10555.0647859239
$ws.Range("C12").Value2 = -381.6682632977227
$ws.Range("G12").Value2 = 60659.45806046302
$ws.Range("H12").Value2 = 3.295288861825552
$ws.Range("B13").Value2 = -10843.62591971395
$ws.Range("C13").Value2 = -288.5611337900527
$ws.Range("G13").Value2 = 65978.94219197202
$ws.Range("H13").Value2 = 3.781339006911555
$ws.Range("B14").Value2 = -11033.59631904412
$ws.Range("C14").Value2 = -189.9703993301763
$ws.Range("G14").Value2 = 71635.50233027429
$ws.Range("H14").Value2 = 4.172457787069916
$ws.Range("B15").Value2 = -11119.37305176982
$ws.Range("C15").Value2 = -85.77673272569882
$ws.Range("G15").Value2 = 77648.81154606806
$ws.Range("H15").Value2 = 4.491294651569389
$ws.Range("B16").Value2 = -11095.23158589149
$ws.Range("C16").Value2 = 24.14146587832965
$ws.Range("G16").Value2 = 84039.63064245693
$ws.Range("H16").Value2 = 4.75403684827842
$ws.Range("B17").Value2 = -10955.32348849623
$ws.Range("C17").Value2 = 139.9080973952632
$ws.Range("G17").Value2 = 90829.8801908141
$ws.Range("H17").Value2 = 4.972536327405597
$ws.Range("B18").Value2 = -10693.6740967304
$ws.Range("C18").Value2 = 261.6493917658336
$ws.Range("G18").Value2 = 98042.71770082673
$ws.Range("H18").Value2 = 5.1556524844107
$ws.Range("B19").Value2 = -10304.18016128637
$ws.Range("C19").Value2 = 389.493935444023
$ws.Range("G19").Value2 = 105702.6202996174
$ws.Range("H19").Value2 = 5.310127733379844
$ws.Range("B20").Value2 = -9780.607462945929
$ws.Range("C20").Value2 = 523.5726983404439
$ws.Range("G20").Value2 = 113835.4733224032
$ws.Range("H20").Value2 = 5.441174624457368
$ws.Range("B21").Value2 = -9116.588402785315
$ws.Range("C21").Value2 = 664.0190601606137
$ws.Range("G21").Value2 = 122468.6652467517
$ws.Range("H21").Value2 = 5.552879246791598
$ws.Range("B22").Value2 = -8305.619566713654
$ws.Range("C22").Value2 = 810.9688360716609
$ws.Range("G22").Value2 = 131631.1894342678
$ws.Range("H22").Value2 = 5.648484409863364
$ws.Range("B23").Value2 = -7341.059265085951
$ws.Range("C23").Value2 = 964.5603016277037
$ws.Range("G23").Value2 = 141353.7531776628
$ws.Range("H23").Value2 = 5.730592247577526
$ws.Range("B24").Value2 = -6216.125048206222
$ws.Range("C24").Value2 = 1124.934216879728
$ws.Range("G24").Value2 = 151668.8945877859
$ws.Range("H24").Value2 = 5.80131165526776
$ws.Range("B25").Value2 = -4923.891198613792
$ws.Range("C25").Value2 = 1292.23384959243
$ws.Range("G25").Value2 = 162611.1078945278
$ws.Range("H25").Value2 = 5.862367230714871
$ws.Range("B26").Value2 = -3457.286201128296
$ws.Range("C26").Value2 = 1466.604997485495
$ws.Range("G26").Value2 = 174216.9777777268
$ws.Range("H26").Value2 = 5.915180888017191
$ws.Range("B27").Value2 = -1809.090191715356
$ws.Range("C27").Value2 = 1648.196009412941
$ws.Range("G27").Value2 = 186525.3233895447
$ws.Range("H27").Value2 = 5.960933769460341
$ws.Range("B28").Value2 = 28.06761367374929
$ws.Range("C28").Value2 = 1837.157805389105
$ws.Range("G28").Value2 = 199577.3527784532
$ws.Range("H28").Value2 = 6.00061375104548
$ws.Range("B29").Value2 = 2061.711509038969
$ws.Range("C29").Value2 = 2033.64389536522
$ws.Range("G29").Value2 = 213416.8284772265
$ws.Range("H29").Value2 = 6.035052277284403
$ws.Range("B30").Value2 = 4299.521905694288
$ws.Range("C30").Value2 = 2237.810396655319
$ws.Range("G30").Value2 = 228090.2450734483
$ws.Range("H30").Value2 = 6.064953198237122
$ws.Range("B31").Value2 = 6749.337955599123
$ws.Range("C31").Value2 = 2449.816049904835
$ws.Range("G31").Value2 = 243647.0196412762
$ws.Range("H31").Value2 = 6.090915546620179
